$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 121, shifting existing rows 121-142 down to 122-143.
$ws.Rows("121:121").Insert()

# Populate the newly inserted row 121 with the new data record.
$ws.Range("A121").Value = 10
$ws.Range("B121").Value = "Vega Modelo de Temuco"
$ws.Range("C121").Value = "La Araucanía"
$ws.Range("D121").Value = 45258
$ws.Range("E121").Value = 9
$ws.Range("F121").Value = 300000001
$ws.Range("G121").Value = "Rabanito"
$ws.Range("H121").Value = "Sin especificar"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 55
$ws.Range("K121").Value = 9000
$ws.Range("L121").Value = 9000
$ws.Range("M121").Value = 9000
$ws.Range("N121").Value = "$/docena de paquetes"
$ws.Range("O121").Value = "Provincia de Cautín"
$ws.Range("P121").Value = 750
$ws.Range("Q121").Value = 12
$ws.Range("R121").Value = "Hortaliza"
